$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 240, shifting existing rows 240:332 down to 241:333
$ws.Rows.Item(240).Insert()

# Populate the newly inserted row 240 with the new data record
$ws.Cells.Item(240, 1).Value = 6
$ws.Cells.Item(240, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(240, 3).Value = "Metropolitana"
$ws.Cells.Item(240, 4).Value = 45120
$ws.Cells.Item(240, 5).Value = 13
$ws.Cells.Item(240, 6).Value = 100112001
$ws.Cells.Item(240, 7).Value = "Berenjena"
$ws.Cells.Item(240, 8).Value = "Sin especificar"
$ws.Cells.Item(240, 9).Value = "Primera"
$ws.Cells.Item(240, 10).Value = 250
$ws.Cells.Item(240, 11).Value = 6000
$ws.Cells.Item(240, 12).Value = 6000
$ws.Cells.Item(240, 13).Value = 6000
$ws.Cells.Item(240, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(240, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(240, 16).Value = 120
$ws.Cells.Item(240, 17).Value = 50
$ws.Cells.Item(240, 18).Value = "Hortaliza"
